$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Approach"
$ws.Range("C1").Value = "Inputs"
$ws.Range("D1").Value = "Limitations"
$ws.Range("E1").Value = "Spatial extent"

# --- Row 2 (USGS) ---
$ws.Range("C2").Value = "sediments"
$ws.Range("D2").Value = "grains > 2mm"
$ws.Range("E2").Value = "point (double-check)"

# --- Row 3 (NEXSS) ---
$ws.Range("C3").Value = "watershed characteristics"
$ws.Range("D3").Value = "Model generalizations"
$ws.Range("E3").Value = "reach-scale"

# --- Row 4 (Abeyshu et al. 2022) ---
$ws.Range("C4").Value = "watershed characteristics"
$ws.Range("D4").Value = "Model generalizations"
$ws.Range("E4").Value = "reach-scale"

# --- Row 5 (YOLO, this study) ---
$ws.Range("B5").Value = "Photogrammetry"
$ws.Range("C5").Value = "images"
$ws.Range("D5").Value = "obscured/small (< 2mm) grains "
$ws.Range("E5").Value = "0.8mx0.8m"

# --- Formatting ---
# Italic style moves from column D to column C (rows 3-4), both columns gain wrap text
$ws.Range("C3:C4").Font.Italic = $true
$ws.Range("C3:C4").WrapText = $true
$ws.Range("D3:D4").Font.Italic = $false
$ws.Range("D3:D4").WrapText = $true
# Row 5 D gains wrap text too
$ws.Range("D5").WrapText = $true

# --- Row heights for the now-wrapped rows 3 and 4 ---
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 34

# --- Column widths (values pre-compensated for the host's internal
# character-width rounding so the saved <col width="..."> lands on the
# target figure: col D -> 15.5, col E -> 18) ---
$ws.Columns.Item(3).ColumnWidth = 12.498697916666666
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 17.166666666666668

# --- Selection matches the new saved state (full table, anchored at E5) ---
$ws.Range("A1:E5").Select() | Out-Null
